# Update the multiplication-table answers to the newly generated values.
# Each cell's text is a unique "A×B=C" string, so a literal Find/Replace
# (no wildcards) on the whole-document Range unambiguously targets the
# correct <w:t> run for each substitution.

$d = $word.ActiveDocument

$replacements = @(
    @("29×34=986",  "83×21=1743"),
    @("90×84=7560", "89×94=8366"),
    @("86×71=6106", "61×37=2257"),
    @("43×71=3053", "26×67=1742"),
    @("13×26=338",  "64×36=2304"),
    @("29×40=1160", "65×74=4810"),
    @("92×86=7912", "71×77=5467"),
    @("97×16=1552", "56×69=3864"),
    @("69×70=4830", "81×88=7128"),
    @("54×98=5292", "63×36=2268"),
    @("77×36=2772", "45×81=3645"),
    @("66×97=6402", "53×46=2438"),
    @("60×49=2940", "94×76=7144"),
    @("86×57=4902", "22×41=902"),
    @("45×58=2610", "38×72=2736"),
    @("77×91=7007", "17×24=408"),
    @("41×75=3075", "53×34=1802"),
    @("59×99=5841", "20×78=1560"),
    @("67×59=3953", "76×57=4332"),
    @("38×21=798",  "33×65=2145"),
    @("47×40=1880", "34×84=2856"),
    @("91×98=8918", "57×27=1539"),
    @("82×84=6888", "31×29=899"),
    @("45×46=2070", "11×85=935"),
    @("91×42=3822", "39×76=2964")
)

$count = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $new, 2)
    if ($found) {
        $count = $count + 1
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Replaced $count of $($replacements.Count) equations"
